$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "Can you identify the elevation of stoery 02?"
$ws.Range("D5").Value = "Can you identify the elevation of level 3?"

$ws.Range("D6").Value = "Tell me the elevation of level 3."

# New strings are introduced in this order so the shared string table is
# rebuilt in the same append sequence as the target workbook.
$ws.Range("C7").Value = "Elevation of floor 03 TO Third Floor."
$ws.Range("C6").Value = "Tell me the elevation of 03 TO Third Floor."
$ws.Range("D7").Value = "Elevation of 02 FL 02 T.O. SLAB."
$ws.Range("D16").Value = "What is the height of the door:300135?"
$ws.Range("E7").Value = "Elevation of 01-EPG."
$ws.Range("E14").Value = "Width of the door:1012722."

$ws.Range("C9").Select()
